# The workbook tracks weekly "Poroto granado" price observations for the
# "Macroferia Regional de Talca" market. This edit adds one more weekly
# observation, inserted as a new row 64 (pushing the former rows 64-83
# down to 65-84, and growing the used range from A1:R83 to A1:R84).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 64; this shifts rows 64:83 down
# to 65:84 and extends the sheet dimension to A1:R84.
$ws.Rows.Item(64).Insert()

# Fill in the data for the newly inserted row 64.
$ws.Cells.Item(64, 1).Value = 5
$ws.Cells.Item(64, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(64, 3).Value = "Maule"
$ws.Cells.Item(64, 4).Value = 44559
$ws.Cells.Item(64, 5).Value = 7
$ws.Cells.Item(64, 6).Value = 100112030
$ws.Cells.Item(64, 7).Value = "Poroto granado"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 300
$ws.Cells.Item(64, 11).Value = 23000
$ws.Cells.Item(64, 12).Value = 23000
$ws.Cells.Item(64, 13).Value = 23000
$ws.Cells.Item(64, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(64, 15).Value = "Región del Maule"
$ws.Cells.Item(64, 16).Value = 920
$ws.Cells.Item(64, 17).Value = 25
$ws.Cells.Item(64, 18).Value = "Hortaliza"
